$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every data cell as a literal inline string (coin name,
# URL, price, % change, date, hour) -- none of it are real Excel numbers.
# For the Price/Volume columns the new text looks like a number or a
# percentage ("257.00", "-1.18%", ...), so a plain .Value assignment would
# get reinterpreted by Excel and silently converted to a numeric cell
# (dropping things like trailing zeros and the literal "%"). Prefixing the
# string with an apostrophe forces Excel to keep it as literal text, and
# resetting the cell Style afterwards clears the "quote prefix" flag that
# the apostrophe entry leaves behind so the cell format matches its
# untouched neighbours.

$ws.Range("D2").Value = "'257.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.18%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.63%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.548"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-12.80%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05892"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.72%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.62%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8576"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.74%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9304"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-7.13%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1410"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.12%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.03589"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.98%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'-2.19%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03231"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.05%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09216"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.32%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001538"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.10%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006037"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-94.35%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006116"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'5.77%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.516"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.59%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.196"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.55%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.60%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3059"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.57%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-1.05%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.856"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.60%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04223"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.35%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001221"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.41%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004291"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-6.22%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.16%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-21.96%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03838"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.67%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006228"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'14.08%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1100"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.95%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002199"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-7.50%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01143"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'5.23%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'0.79%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.15%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.1399"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'63.95%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.1042"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'4,779.06%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.15%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.15%"
$ws.Range("E50").Style = "Normal"
